$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark. In the source document it
#    sits on the empty paragraph right after "U tube manometers"; the
#    edit relocates it, so drop it here first (by name, since it is a
#    hidden bookmark not surfaced through enumeration).
# ------------------------------------------------------------------
$hadBookmark = $false
if ($d.Bookmarks.Exists("_GoBack")) {
    $hadBookmark = $true
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Delete the whole "For singular head losses ... coefficient."
#    paragraph (formula + trailing text), paragraph mark included, so
#    the following paragraph ("The flow regime is defined ...")
#    directly follows the Colebrook-formula paragraph.
# ------------------------------------------------------------------
$singularPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("For s")) {
        $singularPara = $p
        break
    }
}
if ($singularPara -ne $null) {
    $singularPara.Range.Delete()
}

# ------------------------------------------------------------------
# 3) Re-insert the "_GoBack" bookmark, now collapsed at the very start
#    of the "The flow regime is defined ..." paragraph.
# ------------------------------------------------------------------
$flowPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("The flow regime is defined")) {
        $flowPara = $p
        break
    }
}
if ($flowPara -ne $null) {
    $startRange = $flowPara.Range.Duplicate
    $startRange.Collapse(1)
    $d.Bookmarks.Add("_GoBack", $startRange)
}

# ------------------------------------------------------------------
# 4) Drop the stale rendered-page-break marker cached on the
#    "Material" heading by touching the run (self find & replace on
#    just that paragraph's range), which forces the text run to be
#    rewritten without the marker.
# ------------------------------------------------------------------
$materialPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Material") {
        $materialPara = $p
        break
    }
}
if ($materialPara -ne $null) {
    $materialPara.Range.Find.Execute("Material", $true, $false, $false, $false, $false, $true, 1, $false, "Material", 2)
}
